# Clean up the Vancouver Away Defensive Actions sheet:
#  - Unmerge the grouped header cells in row 1
#  - Replace row 1 with clean, per-column header labels
#  - Keep the old (raw) header labels in row 2, but hide that row
#  - Add a hidden, empty spacer row 3
#  - Fill in the missing "Tkl%" (column O) zeros for rows that were missing them
#  - Hide the totals row (row 20)
#  - Nudge the selected cell like the source workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: unmerge the grouped header cells first, so each column can get its own label ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- Row 1: set the clean per-column headers ---
$ws.Cells.Item(1, 1).Value = "Player ID"    # A1
$ws.Cells.Item(1, 2).Value = "Player"       # B1
$ws.Cells.Item(1, 3).Value = "#"            # C1
$ws.Cells.Item(1, 4).Value = "Nation"       # D1
$ws.Cells.Item(1, 5).Value = "Pos"          # E1
$ws.Cells.Item(1, 6).Value = "Age"          # F1
$ws.Cells.Item(1, 7).Value = "90s"          # G1
$ws.Cells.Item(1, 8).Value = "Tkl"          # H1
$ws.Cells.Item(1, 9).Value = "TklW"         # I1
$ws.Cells.Item(1, 10).Value = "Def 3rd"     # J1
$ws.Cells.Item(1, 11).Value = "Mid 3rd"     # K1
$ws.Cells.Item(1, 12).Value = "Att 3rd"     # L1
$ws.Cells.Item(1, 13).Value = "Cha"         # M1
$ws.Cells.Item(1, 14).Value = "Att"         # N1
$ws.Cells.Item(1, 15).Value = "Tkl%"        # O1
$ws.Cells.Item(1, 16).Value = "Lost"        # P1
$ws.Cells.Item(1, 17).Value = "Blocks"      # Q1
$ws.Cells.Item(1, 18).Value = "Sh"          # R1
$ws.Cells.Item(1, 19).Value = "Pass"        # S1
$ws.Cells.Item(1, 20).Value = "Int"         # T1
$ws.Cells.Item(1, 21).Value = "Tkl+Int"     # U1
$ws.Cells.Item(1, 22).Value = "Clr"         # V1
$ws.Cells.Item(1, 23).Value = "Err"         # W1

# --- Row 2 keeps its old raw labels (unchanged) but becomes a hidden helper row ---
$ws.Rows.Item(2).Hidden = $true

# --- Row 3 is an empty hidden spacer row ---
$ws.Rows.Item(3).Hidden = $true

# --- Fill the previously-blank Tkl% (column O) cells with 0 for the data rows that lacked it ---
$oRowsToZero = @(4, 5, 7, 11, 13, 15, 17, 18, 19)
foreach ($r in $oRowsToZero) {
    $ws.Cells.Item($r, 15).Value = 0
}

# --- Hide the totals row ---
$ws.Rows.Item(20).Hidden = $true

# --- Match the saved selection from the source file ---
$ws.Range("O21").Select() | Out-Null
